$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Add the new "link" header
$ws.Range("E1").Value = "link"

# Row 2 becomes RKFox (previously row 3); "inicial" stays 800000000, posicao 1
$ws.Range("A2").Value = "61e852b4dc27dc001969efa3"
$ws.Range("B2").Value = "RKFox"
$ws.Range("C2").Value = 800000000
$ws.Range("D2").Value = 1

# Row 3 becomes HANTAROGAMER (previously row 2); "inicial" changes to 650000000, posicao 2
$ws.Range("A3").Value = "61e484ca5aa1be001868f065"
$ws.Range("B3").Value = "HANTAROGAMER"
$ws.Range("C3").Value = 650000000
$ws.Range("D3").Value = 2

# Fill the new "link" column values in the order the strings were first
# introduced into the shared-string table (row3, row4, then row2).
$ws.Range("E3").Value = "https://rollercoin.com/p/HANTAROGAMER/games"
$ws.Range("E4").Value = "https://rollercoin.com/p/GUERDE/games"
$ws.Range("E2").Value = "https://rollercoin.com/p/RKFox/games"

# Update the selection shown when the workbook is opened
$ws.Range("A5").Select() | Out-Null

# Extend the AutoFilter to the new column and re-sort ascending by "posicao"
$ws.AutoFilterMode = $false
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D1"), $null, 1, $null, $null)
$ws.Sort.SetRange($ws.Range("A1:E4"))
$ws.Sort.Header = 1
$ws.Sort.Apply()
$ws.Range("A1:E1").AutoFilter(1) | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new filter range
$wb.Names.Item("Planilha1!_FilterDatabase").RefersTo = "=Planilha1!`$A`$1:`$E`$1"
